$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.82367267232602
$ws.Range("C2").Value = 0.324682196980101
$ws.Range("B3").Value = 0.175621856814685
$ws.Range("C3").Value = 0.13578283802204

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.01754446692552
$ws.Range("C2").Value = 0.317384281227654
$ws.Range("B3").Value = -0.984930211873869
$ws.Range("C3").Value = 0.101801443691546

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.1068711084345
$ws.Range("C2").Value = 0.179161935426299
$ws.Range("B3").Value = 0.523639577340592
$ws.Range("C3").Value = 0.11890355424349

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.67178985962201
$ws.Range("C2").Value = 0.230673825428547
$ws.Range("B3").Value = 0.0177692748011179
$ws.Range("C3").Value = 0.0167409973495209

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.105418529035825
$ws.Range("B2").Value = -0.0377651881515243
$ws.Range("A3").Value = -0.0377651881515243
$ws.Range("B3").Value = 0.0184369791013194

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.100732781970394
$ws.Range("B2").Value = -0.0284881700528722
$ws.Range("A3").Value = -0.0284881700528722
$ws.Range("B3").Value = 0.010363533937683

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0320989991056973
$ws.Range("B2").Value = -0.00954686235042087
$ws.Range("A3").Value = -0.00954686235042087
$ws.Range("B3").Value = 0.0141380552117345

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0532104137378396
$ws.Range("B2").Value = -0.00252369250313782
$ws.Range("A3").Value = -0.00252369250313782
$ws.Range("B3").Value = 0.000280260992256664
